$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "SelfIntro" test case (row 9) to "ViewSelfIntro", and add two new
# test cases (rows 10, 11) for viewing Hobbies and CCAs.
# Set the "Test Name" (column B) cells first for all three rows.
$ws.Range("B9").Value = "test_<ViewSelfIntro>"
$ws.Range("B10").Value = "test_<ViewHobbies>"
$ws.Range("B11").Value = "test_<ViewCCA>"

# Row 9: test description text is unchanged.
$ws.Range("C9").Value = "This is to test whether users are able to view a created SelfIntro"

# Test descriptions for the two new rows.
$ws.Range("C10").Value = "This is to test whether users are able to view a list of hobbies"
$ws.Range("C11").Value = "This is to test whether users are able to view a list of CCAs"

# Test values (NIL) for all three rows.
$ws.Range("D9").Value = "NIL"
$ws.Range("D10").Value = "NIL"
$ws.Range("D11").Value = "NIL"

# Expected outcomes.
$ws.Range("E9").Value = "Self Intro page is shown"
$ws.Range("E10").Value = "Hobbies are shown"
$ws.Range("E11").Value = "CCAs are shown"

# Outcome for row 9.
$ws.Range("F9").Value = "Self Intro page is shown"

# New row numbering for the two added test cases.
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Update the active selection to match the saved view state
$ws.Range("F15").Select()
